# Refresh the crypto symbol list (GitHub Actions scheduled update).
# All Price/Volume(1h) cells are stored as text in the sheet, so numeric-looking
# values are written with a leading apostrophe to keep them as text (matches the
# existing inline-string cell type instead of letting Excel auto-convert them to
# numbers/percentages).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.79"
$ws.Range("E2").Value = "'0.88%"
$ws.Range("D3").Value = "'29.13"
$ws.Range("E3").Value = "'-2.44%"
$ws.Range("D4").Value = "'5.143"
$ws.Range("E4").Value = "'0.01%"
$ws.Range("E5").Value = "'1.99%"
$ws.Range("D6").Value = "'6.612"
$ws.Range("E6").Value = "'1.38%"
$ws.Range("D7").Value = "'3.175"
$ws.Range("E7").Value = "'5.07%"
$ws.Range("D8").Value = "'0.8568"
$ws.Range("E8").Value = "'1.63%"
$ws.Range("D9").Value = "'0.8619"
$ws.Range("E9").Value = "'0.44%"
$ws.Range("D10").Value = "'0.1375"
$ws.Range("E10").Value = "'2.84%"
$ws.Range("D11").Value = "'0.07077"
$ws.Range("E11").Value = "'2.45%"
$ws.Range("D12").Value = "'0.03290"
$ws.Range("E12").Value = "'13.92%"
$ws.Range("D13").Value = "'0.09359"
$ws.Range("E13").Value = "'-0.30%"
$ws.Range("D14").Value = "'0.001522"
$ws.Range("E14").Value = "'-0.24%"
$ws.Range("D15").Value = "'0.0005993"
$ws.Range("E15").Value = "'-94.06%"
$ws.Range("D16").Value = "'0.005975"
$ws.Range("E16").Value = "'-1.98%"
$ws.Range("D17").Value = "'3.487"
$ws.Range("E17").Value = "'-0.60%"
$ws.Range("D18").Value = "'2.198"
$ws.Range("E18").Value = "'-1.98%"
$ws.Range("E19").Value = "'1.58%"
$ws.Range("D20").Value = "'0.03333"
$ws.Range("E20").Value = "'1.91%"
$ws.Range("D21").Value = "'0.1281"
$ws.Range("E21").Value = "'-1.65%"
$ws.Range("D22").Value = "'3.176"
$ws.Range("E22").Value = "'-11.94%"
$ws.Range("D23").Value = "'0.04136"
$ws.Range("E23").Value = "'-0.83%"
$ws.Range("D24").Value = "'0.1398"
$ws.Range("E24").Value = "'1.87%"
$ws.Range("E25").Value = "'1.28%"
$ws.Range("D26").Value = "'0.004142"
$ws.Range("E26").Value = "'-6.74%"
$ws.Range("E27").Value = "'2.55%"
$ws.Range("E28").Value = "'3.38%"
$ws.Range("D40").Value = "'0.03729"
$ws.Range("E40").Value = "'0.48%"
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D41").Value = "'0.005764"
$ws.Range("E41").Value = "'8.26%"
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").Value = "'0.1069"
$ws.Range("E42").Value = "'1.23%"
$ws.Range("D43").Value = "'0.002198"
$ws.Range("E43").Value = "'-3.89%"
$ws.Range("D44").Value = "'0.009168"
$ws.Range("D45").Value = "'0.00005283"
$ws.Range("E45").Value = "'3.79%"
$ws.Range("E46").Value = "'0.01%"
$ws.Range("D47").Value = "'0.05793"
$ws.Range("E47").Value = "'-42.00%"
$ws.Range("D48").Value = "'0.002172"
$ws.Range("E48").Value = "'-22.65%"
$ws.Range("E49").Value = "'0.01%"
$ws.Range("E50").Value = "'0.01%"
